$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("variables")

$ws.Range("A18").Value = "dt"
$ws.Range("B18").Value = 0.00001
$ws.Range("B18").NumberFormat = "0.00E+00"
$ws.Range("C18").Value = "s"
$ws.Range("D18").Value = "time step in runge-kutta integration"

$ws.Range("F13").Select()
